$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes")

# Add the new row of data (row 18) documenting the Custom Fields behavior.
$ws.Range("A18").Value = 6
$ws.Range("B18").Value = "Level 0-0"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = "Modify"
$ws.Range("E18").Value = "WSJF"
$ws.Range("F18").Value = 3.14159
$ws.Range("H18").Value = "Custom Fields on a board are checked if it is not a standard Card field. If still no match, then it is ignored"

# Match formatting of the other note cells in column H (text format + wrap).
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").WrapText = $true

# The new row uses the taller row height seen on other multi-line rows.
$ws.Rows.Item(18).RowHeight = 30

# Make "Changes" the active sheet/tab and select the newly added note cell,
# matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("H18").Select() | Out-Null
